{"js": "// Split the single run of the Title, Author and Abstract paragraphs into\n// one run per word plus one run per inter-word space, leaving each\n// paragraph's visible text unchanged.\n//\n// Office.js has no \"split this run at offset N\" primitive, so rebuild\n// the affected paragraphs' run lists directly via Range.insertOoxml:\n// build a Flat-OPC-wrapped <w:p> (keeping the paragraph's existing\n// style) containing the desired sequence of <w:r><w:t> runs, then\n// replace the paragraph's range with it. That swaps the run list in\n// place without adding or removing paragraphs.\n\nfunction escapeXml(s) {\n  return s.replace(/&/g, \"&amp;\").replace(/</g, \"&lt;\").replace(/>/g, \"&gt;\");\n}\n\n// \"A selection of questions\" -> [\"A\", \" \", \"selection\", \" \", \"of\", \" \", \"questions\"]\n// Split on the space character, keeping each space as its own segment;\n// everything else (including punctuation) stays glued to its word.\nfunction splitIntoWordRuns(text) {\n  const parts = [];\n  let word = \"\";\n  for (const ch of text) {\n    if (ch === \" \") {\n      if (word.length > 0) {\n        parts.push(word);\n        word = \"\";\n      }\n      parts.push(\" \");\n    } else {\n      word += ch;\n    }\n  }\n  if (word.length > 0) parts.push(word);\n  return parts;\n}\n\nfunction buildFlatOpcParagraph(styleId, parts) {\n  const runsXml = parts\n    .map((p) => `<w:r><w:t xml:space=\"preserve\">${escapeXml(p)}</w:t></w:r>`)\n    .join(\"\");\n  const pPr = styleId ? `<w:pPr><w:pStyle w:val=\"${styleId}\"/></w:pPr>` : \"\";\n  return (\n    `<?xml version=\"1.0\"?>` +\n    `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n    `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">` +\n    `<pkg:xmlData>` +\n    `<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">` +\n    `<w:body><w:p>${pPr}${runsXml}</w:p></w:body>` +\n    `</w:document>` +\n    `</pkg:xmlData></pkg:part></pkg:package>`\n  );\n}\n\nconst targetStyles = [\"Title\", \"Author\", \"Abstract\"];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  paragraph.load(\"text,style\");\n}\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (targetStyles.includes(paragraph.style)) {\n    const parts = splitIntoWordRuns(paragraph.text);\n    const flatOpcXml = buildFlatOpcParagraph(paragraph.style, parts);\n    const range = paragraph.getRange();\n    range.insertOoxml(flatOpcXml, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Split the single w:r run of a few title/byline/abstract paragraphs into\n# one w:r per word plus one w:r per inter-word space, leaving the visible\n# text of each paragraph unchanged.\n#\n# Word's COM object model has no direct \"split this run at offset N\"\n# call, so we rebuild each target paragraph's runs via a Flat-OPC\n# `Range.InsertXML` replace: construct a <w:p> with the paragraph's\n# existing style plus the desired sequence of <w:r><w:t> runs, then\n# insert it over the paragraph's current Range (text only, not the\n# paragraph mark). That swaps the run list in place without adding or\n# removing any paragraphs.\n\nfunction ConvertTo-XmlText($text) {\n    $text.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')\n}\n\nfunction Split-TextIntoWordRuns($text) {\n    # \"A selection of questions\" -> \"A\" \" \" \"selection\" \" \" \"of\" \" \" \"questions\"\n    # i.e. split on the space character, keeping each space as its own\n    # segment, everything else (including punctuation) stays glued to\n    # its neighbouring word.\n    $parts = New-Object System.Collections.Generic.List[string]\n    $word = \"\"\n    foreach ($ch in $text.ToCharArray()) {\n        if ($ch -eq ' ') {\n            if ($word.Length -gt 0) {\n                $parts.Add($word)\n                $word = \"\"\n            }\n            $parts.Add(\" \")\n        } else {\n            $word += $ch\n        }\n    }\n    if ($word.Length -gt 0) {\n        $parts.Add($word)\n    }\n    return $parts\n}\n\nfunction New-FlatOpcParagraphXml($styleName, $parts) {\n    $runsXml = \"\"\n    foreach ($part in $parts) {\n        $runsXml += '<w:r><w:t xml:space=\"preserve\">' + (ConvertTo-XmlText $part) + '</w:t></w:r>'\n    }\n    $pPr = \"\"\n    if ($styleName) {\n        $pPr = '<w:pPr><w:pStyle w:val=\"' + $styleName + '\"/></w:pPr>'\n    }\n    return '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $pPr + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\nfunction Expand-ParagraphRuns($paragraph) {\n    $styleName = $paragraph.Style.NameLocal\n    $rng = $paragraph.Range\n    # Paragraph.Range.Text carries a trailing \"\\r\" for the paragraph mark;\n    # trim it so it doesn't become part of the last run's text.\n    $text = $rng.Text.TrimEnd([char]13)\n    $parts = Split-TextIntoWordRuns $text\n    $rng.InsertXML((New-FlatOpcParagraphXml $styleName $parts))\n}\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    if ($styleName -eq \"Title\" -or $styleName -eq \"Author\" -or $styleName -eq \"Abstract\") {\n        Expand-ParagraphRuns $p\n    }\n}\n"}
